$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.450.80"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +2.85%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.314.53"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.78%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.08"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.85"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +5.22%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.334.22"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.03%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.06%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.12%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.341"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.35%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.91"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.729.32"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +1.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.536.30"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.301.01"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.04%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "323.78"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.57"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.64"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.46%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.995"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.98"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.27"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +10.94%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.06%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.41"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.27%  "
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.70"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.20"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.41"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.46%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.994"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.55%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.28%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +6.94%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.96"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.381"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +1.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.62"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +3.67%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.22"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +6.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "275.33"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +6.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0936"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.05%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.87%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.58%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.48%  "
